$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 5
$ws.Range("G2").Value = 7.8
$ws.Range("H2").Value = 1.52
$ws.Range("I2").Value = 1.66
$ws.Range("J2").Value = 3.85
$ws.Range("K2").Value = 5.2
$ws.Range("M2").Value = 1.03
$ws.Range("N2").Value = 4
$ws.Range("O2").Value = 1.2
$ws.Range("P2").Value = 2.28
$ws.Range("Q2").Value = 1.61
$ws.Range("R2").Value = 1.52
$ws.Range("S2").Value = 2.32
$ws.Range("T2").Value = 1.72
$ws.Range("U2").Value = 2.08
$ws.Range("V2").Value = 2.52
$ws.Range("W2").Value = 1.16
$ws.Range("X2").Value = 24
$ws.Range("Y2").Value = 11.5
$ws.Range("Z2").Value = 11.5
$ws.Range("AA2").Value = 16
$ws.Range("AB2").Value = 27
$ws.Range("AC2").Value = 11.5
$ws.Range("AD2").Value = 11
$ws.Range("AE2").Value = 16.5
$ws.Range("AG2").Value = 26
$ws.Range("AH2").Value = 21
$ws.Range("AO2").Value = 7.2
$ws.Range("K3").Value = 980
$ws.Range("Q3").Value = 1.02
$ws.Range("G4").Value = 1.35
$ws.Range("K4").Value = 5.9
$ws.Range("L4").Value = 1.28
$ws.Range("P4").Value = 2.66
$ws.Range("W4").Value = 3.85
$ws.Range("F5").Value = 5.3
$ws.Range("G5").Value = 5.5
$ws.Range("H5").Value = 1.73
$ws.Range("I5").Value = 1.74
$ws.Range("V5").Value = 2.34
$ws.Range("W5").Value = 1.22
$ws.Range("AJ5").Value = 130
$ws.Range("L7").Value = 1.25
$ws.Range("J8").Value = 3.6
$ws.Range("P8").Value = 2.08
$ws.Range("F9").Value = 1.34
$ws.Range("L9").Value = 1.2
$ws.Range("P9").Value = 3.6
$ws.Range("U9").Value = 2.44
$ws.Range("W9").Value = 3.85
$ws.Range("AA9").Value = 570
$ws.Range("AB9").Value = 16
$ws.Range("AL9").Value = 23
$ws.Range("AO9").Value = 80
$ws.Range("G10").Value = 1.41
$ws.Range("O10").Value = 1.18
$ws.Range("P10").Value = 2.7
$ws.Range("Q10").Value = 1.56
$ws.Range("S10").Value = 2.38
$ws.Range("V10").Value = 1.11
$ws.Range("W10").Value = 3.45
$ws.Range("AC10").Value = 12.5
$ws.Range("AL10").Value = 28
$ws.Range("AN10").Value = 4.7
$ws.Range("H11").Value = 5.3
$ws.Range("I11").Value = 5.4
$ws.Range("Q11").Value = 1.59
$ws.Range("R11").Value = 1.64
$ws.Range("V11").Value = 1.22
$ws.Range("Z11").Value = 46
$ws.Range("AF11").Value = 12
$ws.Range("F12").Value = 2.3
$ws.Range("G12").Value = 2.34
$ws.Range("S12").Value = 2.98
$ws.Range("W12").Value = 1.75
$ws.Range("AA12").Value = 60
$ws.Range("AM12").Value = 70
$ws.Range("AN12").Value = 14
$ws.Range("J13").Value = 5.3
$ws.Range("P13").Value = 2.72
$ws.Range("Q13").Value = 1.55
$ws.Range("R13").Value = 1.71
$ws.Range("AK13").Value = 95
$ws.Range("AM13").Value = 95
$ws.Range("F14").Value = 2.76
$ws.Range("H14").Value = 2.66
$ws.Range("I14").Value = 2.68
$ws.Range("J14").Value = 3.7
$ws.Range("L14").Value = 1.32
$ws.Range("M14").Value = 1.05
$ws.Range("V14").Value = 1.59
$ws.Range("AN14").Value = 17.5
$ws.Range("Q15").Value = 1.02
$ws.Range("F16").Value = 1.91
$ws.Range("G16").Value = 2.04
$ws.Range("H16").Value = 4
$ws.Range("I16").Value = 4.5
$ws.Range("J16").Value = 3.7
$ws.Range("Q16").Value = 1.71
$ws.Range("R16").Value = 1.41
$ws.Range("S16").Value = 3.05
$ws.Range("T16").Value = 1.73
$ws.Range("U16").Value = 2.18
$ws.Range("V16").Value = 1.28
$ws.Range("W16").Value = 1.96
$ws.Range("X16").Value = 1000
$ws.Range("AA16").Value = 100
$ws.Range("AB16").Value = 10.5
$ws.Range("AC16").Value = 9
$ws.Range("AF16").Value = 13.5
$ws.Range("AG16").Value = 10.5
$ws.Range("AM16").Value = 110
$ws.Range("AN16").Value = 980
$ws.Range("N17").Value = 1.29
$ws.Range("P17").Value = 1.28
$ws.Range("F18").Value = 1.78
$ws.Range("G18").Value = 1.95
$ws.Range("H18").Value = 4.6
$ws.Range("I18").Value = 6
$ws.Range("J18").Value = 3.35
$ws.Range("K18").Value = 4.3
$ws.Range("M18").Value = 1.09
$ws.Range("N18").Value = 3
$ws.Range("O18").Value = 1.42
$ws.Range("P18").Value = 1.68
$ws.Range("Q18").Value = 2.22
$ws.Range("R18").Value = 1.25
$ws.Range("S18").Value = 4.2
$ws.Range("T18").Value = 2.02
$ws.Range("U18").Value = 1.81
$ws.Range("V18").Value = 1.2
$ws.Range("W18").Value = 2.04
$ws.Range("AI18").Value = 110
$ws.Range("G19").Value = 2.24
$ws.Range("U19").Value = 1.64
$ws.Range("AI19").Value = 150
$ws.Range("AM19").Value = 310
$ws.Range("Z20").Value = 1000
$ws.Range("AO20").Value = 17
$ws.Range("N21").Value = 3.9
$ws.Range("P21").Value = 2.16
$ws.Range("U21").Value = 2
$ws.Range("AB21").Value = 11.5
$ws.Range("AC21").Value = 1000
$ws.Range("AF21").Value = 12
$ws.Range("AM21").Value = 130
$ws.Range("AN21").Value = 1000
